# corrected refclk placement, LHC clock was incorrectly placed as refclk
# Adds a new "Swatch index" column (I) for the CPPF link-map block (rows 44-51):
#   I44 = header "Swatch index"
#   I45..I51 = 1..7 (numeric swatch index values)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I44").Value = "Swatch index"

$ws.Range("I45").Value = 1
$ws.Range("I46").Value = 2
$ws.Range("I47").Value = 3
$ws.Range("I48").Value = 4
$ws.Range("I49").Value = 5
$ws.Range("I50").Value = 6
$ws.Range("I51").Value = 7

# reflect the author's last-selected cell after making the edit
$ws.Range("I47").Select()
